$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.069.92'
$ws.Range('E2').Value = '  +3.71%  '

$ws.Range('D3').Value = '2.665.56'
$ws.Range('E3').Value = '  +6.23%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.18'
$ws.Range('E5').Value = '  +6.25%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '328.83'
$ws.Range('E6').Value = '  +2.80%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('E7').Value = '  +1.13%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  +3.72%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.43'
$ws.Range('E10').Value = '  +6.06%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.22'
$ws.Range('E11').Value = '  +1.61%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0828'
$ws.Range('E12').Value = '  +2.47%  '

$ws.Range('E13').Value = '  +0.87%  '

$ws.Range('E14').Value = '  +4.09%  '

$ws.Range('D15').Value = '3.021.17'
$ws.Range('E15').Value = '  +4.06%  '

$ws.Range('D16').Value = '2.677.08'
$ws.Range('E16').Value = '  +6.54%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.878'
$ws.Range('E17').Value = '  +5.26%  '

$ws.Range('D18').Value = '50.020.46'
$ws.Range('E18').Value = '  +3.85%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.26'
$ws.Range('E19').Value = '  +1.99%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.81'
$ws.Range('E20').Value = '  +2.33%  '

$ws.Range('E21').Value = '  -0.80%  '

$ws.Range('E22').Value = '  +2.71%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.59'
$ws.Range('E23').Value = '  +1.86%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '278.79'
$ws.Range('E24').Value = '  +2.17%  '

$ws.Range('E25').Value = '  +2.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.99'
$ws.Range('E26').Value = '  +3.93%  '

$ws.Range('E27').Value = '  -0.05%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.99'
$ws.Range('E28').Value = '  +2.42%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.22'
$ws.Range('E29').Value = '  -2.87%  '

$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.41'
$ws.Range('E30').Value = '  +4.87%  '

$ws.Range('E31').Value = '  -2.49%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.52'
$ws.Range('E32').Value = '  +1.92%  '

$ws.Range('E33').Value = '  +2.44%  '

$ws.Range('E34').Value = '  +1.73%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0802'
$ws.Range('E35').Value = '  +2.62%  '

$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.09'
$ws.Range('E37').Value = '  +7.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.79'
$ws.Range('E38').Value = '  +3.63%  '

$ws.Range('E39').Value = '  +6.91%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '125.98'
$ws.Range('E40').Value = '  +4.82%  '

$ws.Range('E41').Value = '  +1.79%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.74'
$ws.Range('E42').Value = '  +3.64%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.24'
$ws.Range('E43').Value = '  +1.06%  '

$ws.Range('E44').Value = '  +3.60%  '

$ws.Range('E45').Value = '  +5.49%  '

$ws.Range('D46').Value = '2.077.28'
$ws.Range('E46').Value = '  +3.59%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.25'
$ws.Range('E47').Value = '  +12.18%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.00'
$ws.Range('E48').Value = '  +5.71%  '

$ws.Range('E49').Value = '  +2.10%  '

$ws.Range('E50').Value = '  +3.34%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '82.24'
$ws.Range('E51').Value = '  +3.51%  '

$ws.Range('B2:E51').ClearFormats()
